# update template skeletons - subtitles, authors
#
# The "Hyperlink" character style used to render as italic accent1-blue;
# it now renders as non-italic, theme "Text 1" (near-black).
$d = $word.ActiveDocument

$hyperlink = $d.Styles("Hyperlink")

# <w:i/>  ->  <w:i w:val="0"/>
$hyperlink.Font.Italic = $false

# <w:color w:val="4F81BD" w:themeColor="accent1"/>
#   ->
# <w:color w:val="000000" w:themeColor="text1"/>
$hyperlink.Font.TextColor.ObjectThemeColor = 13
